$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings that must stay as plain text
# (matching the source data), so force text format before assigning each one.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.797.94"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.635.89"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.32"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0643"
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("E10").Value = "  +1.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0784"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.640.07"
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.861.31"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("E15").Value = "  -1.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₃0777"
$ws.Range("E16").Value = "  +2.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.08"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.818.62"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.45"
$ws.Range("E20").Value = "  +2.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.99"
$ws.Range("E21").Value = "  -0.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.95"
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("E23").Value = "  +0.57%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "139.25"
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("E27").Value = "  -4.72%  "
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.56"
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("E31").Value = "  +1.89%  "
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.60"
$ws.Range("E34").Value = "  +2.31%  "
$ws.Range("E35").Value = "  +0.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.898"
$ws.Range("E36").Value = "  -0.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.57"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.552"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.106.60"
$ws.Range("E39").Value = "  -2.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0157"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.59"
$ws.Range("E42").Value = "  +0.85%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.27"
$ws.Range("E43").Value = "  +1.59%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.802"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₆0110"
$ws.Range("E45").Value = "  -2.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.52"
$ws.Range("E47").Value = "  +12.85%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.73"
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.418"
$ws.Range("E49").Value = "  -5.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0504"
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("E51").Value = "  -0.01%  "
